# Fix typo in generated file names (TLC01 -> TCL01) and update the
# sheet's saved view/selection to match the latest manual review state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "TC01_ICDC_TCL01_Breed-GoldenRetriever_TSVData.xlsx"
$ws.Range("E2").Value = "TC01_ICDC_TCL01_Breed-GoldenRetriever_WebData.xlsx"

$ws.Range("D2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 1
